$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.190.65"
$ws.Range("E2").Value = "  -4.56%  "

# Row 3
$ws.Range("D3").Value = "1.756.28"
$ws.Range("E3").Value = "  -3.53%  "

# Row 4
$ws.Range("D4").Value = "0.9830"
$ws.Range("E4").Value = "  -2.00%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "304.69"
$ws.Range("E5").Value = "  -2.41%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "0.9901"
$ws.Range("E6").Value = "  -1.22%  "

# Row 7
$ws.Range("D7").Value = "0.4250"
$ws.Range("E7").Value = "  +0.39%  "

# Row 8
$ws.Range("D8").Value = "0.3626"
$ws.Range("E8").Value = "  +1.04%  "

# Row 9
$ws.Range("D9").Value = "0.07134"
$ws.Range("E9").Value = "  -0.50%  "

# Row 10
$ws.Range("D10").Value = "0.8438"
$ws.Range("E10").Value = "  -0.98%  "

# Row 11
$ws.Range("D11").Value = "20.26"
$ws.Range("E11").Value = "  -0.58%  "

# Row 12
$ws.Range("D12").Value = "1.782.65"
$ws.Range("E12").Value = "  -3.92%  "

# Row 13
$ws.Range("D13").Value = "6.394"
$ws.Range("E13").Value = "  -0.34%  "

# Row 14
$ws.Range("D14").Value = "5.211"
$ws.Range("E14").Value = "  -2.51%  "

# Row 15
$ws.Range("D15").Value = "0.06833"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16
$ws.Range("D16").Value = "0.9952"
$ws.Range("E16").Value = "  -1.02%  "

# Row 17
$ws.Range("D17").Value = "78.33"
$ws.Range("E17").Value = "  -3.50%  "

# Row 18
$ws.Range("D18").Value = "0.000008629"
$ws.Range("E18").Value = "  -2.64%  "

# Row 19
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  -0.28%  "

# Row 20
$ws.Range("D20").Value = "14.92"
$ws.Range("E20").Value = "  -2.36%  "

# Row 21
$ws.Range("D21").Value = "26.273.61"
$ws.Range("E21").Value = "  -4.03%  "

# Row 22
$ws.Range("D22").Value = "5.065"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").Value = "2.002.54"
$ws.Range("E24").Value = "  -3.88%  "

# Row 25
$ws.Range("D25").Value = "150.79"
$ws.Range("E25").Value = "  -2.18%  "

# Row 26
$ws.Range("D26").Value = "1.844"
$ws.Range("E26").Value = "  -7.04%  "

# Row 27
$ws.Range("D27").Value = "17.95"
$ws.Range("E27").Value = "  -1.98%  "

# Row 28
$ws.Range("D28").Value = "5.067"
$ws.Range("E28").Value = "  -0.93%  "

# Row 29
$ws.Range("D29").Value = "113.00"
$ws.Range("E29").Value = "  -0.81%  "

# Row 30
$ws.Range("D30").Value = "1.780"
$ws.Range("E30").Value = "  +0.89%  "

# Row 31
$ws.Range("D31").Value = "0.08891"
$ws.Range("E31").Value = "  -0.09%  "

# Row 32
$ws.Range("D32").Value = "0.7237"
$ws.Range("E32").Value = "  -2.68%  "

# Row 33
$ws.Range("D33").Value = "1.109"
$ws.Range("E33").Value = "  -0.36%  "

# Row 34
$ws.Range("D34").Value = "4.294"
$ws.Range("E34").Value = "  -4.79%  "

# Row 35
$ws.Range("D35").Value = "0.9962"
$ws.Range("E35").Value = "  -0.65%  "

# Row 36
$ws.Range("D36").Value = "2.701"
$ws.Range("E36").Value = "  -8.17%  "

# Row 37
$ws.Range("D37").Value = "1.070"
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("D38").Value = "0.05122"
$ws.Range("E38").Value = "  -1.87%  "

# Row 39
$ws.Range("D39").Value = "0.01869"
$ws.Range("E39").Value = "  -2.41%  "

# Row 40
$ws.Range("D40").Value = "0.4884"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41
$ws.Range("D41").Value = "0.1597"
$ws.Range("E41").Value = "  -3.12%  "

# Row 42
$ws.Range("D42").Value = "2.553"
$ws.Range("E42").Value = "  -7.90%  "

# Row 43
$ws.Range("D43").Value = "6.243"
$ws.Range("E43").Value = "  -0.76%  "

# Row 44
$ws.Range("D44").Value = "7.952"
$ws.Range("E44").Value = "  -3.86%  "

# Row 45
$ws.Range("D45").Value = "104.25"
$ws.Range("E45").Value = "  -1.49%  "

# Row 46
$ws.Range("D46").Value = "10.17"
$ws.Range("E46").Value = "  -2.32%  "

# Row 47
$ws.Range("D47").Value = "0.9870"
$ws.Range("E47").Value = "  -1.48%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.4464"
$ws.Range("E48").Value = "  -3.84%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06139"
$ws.Range("E49").Value = "  -4.66%  "

# Row 50
$ws.Range("D50").Value = "1.593"
$ws.Range("E50").Value = "  -0.98%  "

# Row 51
$ws.Range("D51").Value = "1.731"
$ws.Range("E51").Value = "  +3.24%  "
